$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 08:15:36"
$wsZhCn.Range("E3").Value = "2016-03-12 08:15:36"
$wsZhCn.Range("H2").Value = "2016-03-12 08:15:54"
$wsZhCn.Range("H3").Value = "2016-03-12 08:15:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 08:15:39"
$wsDeDe.Range("E3").Value = "2016-03-12 08:15:39"
$wsDeDe.Range("H2").Value = "2016-03-12 08:16:00"
$wsDeDe.Range("H3").Value = "2016-03-12 08:16:00"
